$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$noise = @{
  2 = @(0.75, 0.7333333333333333)
  3 = @(0, 0.1666666666666667)
  4 = @(0.5, 0.6666666666666666)
  5 = @(0.25, 0.4444444444444444)
  6 = @(0.5, 0.5555555555555556)
  7 = @(0.75, 0.4444444444444444)
  8 = @(0.5, 0.7777777777777778)
  9 = @(0.5, 0.3333333333333333)
  10 = @(0.75, 0.7333333333333333)
  11 = @(0.25, 0.4166666666666667)
  12 = @(0.25, 0.5333333333333333)
  13 = @(0, 0.3333333333333333)
  14 = @(0.75, 0.5833333333333334)
  15 = @(0.5, 0.6666666666666666)
  16 = @(0.25, 0.5833333333333334)
  17 = @(0, 0.3333333333333333)
  18 = @(0.5, 0.6666666666666666)
  19 = @(0.75, 0.5)
  20 = @(0.5, 0.6)
  21 = @(0, 0.1666666666666667)
  22 = @(0.5, 0.3333333333333333)
  23 = @(0.25, 0.4166666666666667)
  24 = @(0.5, 0.6666666666666666)
  25 = @(0.5, 0.5)
  26 = @(0.5, 0.6666666666666666)
  27 = @(0.75, 0.3333333333333333)
  28 = @(0.5, 0.6)
  29 = @(0.5, 0.7333333333333333)
  30 = @(0.5, 0.75)
  31 = @(0, 0.3333333333333333)
  32 = @(0.75, 0.6666666666666666)
  33 = @(0, 0.5555555555555556)
  34 = @(0.25, 0.5833333333333334)
  35 = @(0, 0.3333333333333333)
  36 = @(0.5, 0.6666666666666666)
  37 = @(0.5, 0.5)
  38 = @(0.75, 0.7333333333333333)
  39 = @(0, 0.25)
  40 = @(0.5, 0.6666666666666666)
  41 = @(0.75, 0.4444444444444444)
  42 = @(0, 0.3333333333333333)
  43 = @(0.5, 0.5)
  44 = @(0.25, 0.2222222222222222)
  45 = @(0.75, 0.5)
  46 = @(0.5, 0.6666666666666666)
  47 = @(0, 0.3333333333333333)
  48 = @(0.5, 0.5)
  49 = @(0.75, 0.6666666666666666)
  50 = @(0.25, 0.4166666666666667)
  51 = @(0.5, 0.5555555555555556)
  52 = @(0.75, 0.6666666666666666)
  53 = @(0, 0.2)
  54 = @(0.75, 0.4444444444444444)
  55 = @(0.5, 0.4444444444444444)
  56 = @(0, 0.1666666666666667)
  57 = @(0.75, 0.5)
  58 = @(0, 0.1666666666666667)
  59 = @(0.25, 0.4444444444444444)
  60 = @(0, 0.3333333333333333)
  61 = @(0.75, 0.5)
  62 = @(0, 0.1666666666666667)
  63 = @(0.5, 0.6)
  64 = @(0, 0.3333333333333333)
  65 = @(0.75, 0.6666666666666666)
  66 = @(0.75, 0.75)
  67 = @(0.75, 0.4444444444444444)
  68 = @(0, 0.3333333333333333)
  69 = @(0.5, 0.5)
  70 = @(0.75, 0.7333333333333333)
  71 = @(0.75, 0.75)
  72 = @(0, 0.3333333333333333)
  73 = @(0.5, 0.5)
  74 = @(0, 0.3333333333333333)
  75 = @(0.5, 0.6)
  76 = @(0, 0.3333333333333333)
  77 = @(0.5, 0.5555555555555556)
  78 = @(0, 0.3333333333333333)
  79 = @(0.5, 0.6666666666666666)
  80 = @(0.5, 0.6666666666666666)
  81 = @(0.75, 0.6666666666666666)
  82 = @(0, 0.3333333333333333)
  83 = @(0.75, 0.6666666666666666)
  84 = @(0.5, 0.4444444444444444)
  85 = @(0.75, 0.6666666666666666)
  86 = @(0.75, 0.6666666666666666)
  87 = @(0.5, 0.5555555555555556)
  88 = @(0.75, 0.5)
  89 = @(0, 0.1666666666666667)
  90 = @(0.5, 0.6666666666666666)
  91 = @(0.75, 0.4444444444444444)
  92 = @(0, 0.1666666666666667)
  93 = @(0.75, 0.5)
  94 = @(0.25, 0.4444444444444444)
  95 = @(0.5, 0.5555555555555556)
  96 = @(0, 0.3333333333333333)
  97 = @(0.5, 0.6666666666666666)
  98 = @(0.75, 0.5)
  99 = @(0, 0.3333333333333333)
  100 = @(0.5, 0.5333333333333333)
  101 = @(0.5, 0.5555555555555556)
}

foreach ($r in $noise.Keys) {
  $vals = $noise[$r]
  $ws.Cells.Item([int]$r, 2).Value = $vals[0]
  $ws.Cells.Item([int]$r, 3).Value = $vals[1]
}